# Apply the Sun Apr 30 2023 GitHub Actions cryptos-list refresh to the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) and Volume(1h) (column E) refresh for each coin row.
# Cells whose new text would otherwise be auto-parsed as a number by Excel
# (e.g. "1.014", "0.06877") are explicitly formatted as Text ("@") first so
# the literal string is preserved exactly, matching the source data feed.

$ws.Range('D2').Value = '29.610.49'
$ws.Range('E2').Value = '  +0.39%  '

$ws.Range('D3').Value = '1.925.60'
$ws.Range('E3').Value = '  +0.21%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.014'
$ws.Range('E4').Value = '  +0.85%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '326.89'
$ws.Range('E5').Value = '  +0.39%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.012'
$ws.Range('E6').Value = '  +0.78%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4821'
$ws.Range('E7').Value = '  -0.46%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4053'
$ws.Range('E8').Value = '  -0.87%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08189'
$ws.Range('E9').Value = '  +0.15%  '

$ws.Range('E10').Value = '  -1.75%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '23.69'
$ws.Range('E11').Value = '  -0.40%  '

$ws.Range('D12').Value = '1.905.78'
$ws.Range('E12').Value = '  +0.94%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.059'
$ws.Range('E13').Value = '  +0.18%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.280'
$ws.Range('E14').Value = '  +0.54%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.39'
$ws.Range('E15').Value = '  -0.26%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.06877'
$ws.Range('E16').Value = '  +1.83%  '

$ws.Range('E17').Value = '  +0.70%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001038'
$ws.Range('E18').Value = '  -0.17%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.60'
$ws.Range('E19').Value = '  -0.89%  '

$ws.Range('E20').Value = '  +0.57%  '

$ws.Range('D21').Value = '29.595.68'
$ws.Range('E21').Value = '  +0.25%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.644'
$ws.Range('E22').Value = '  +0.15%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.96'
$ws.Range('E23').Value = '  +1.62%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.204'
$ws.Range('E24').Value = '  +1.10%  '

$ws.Range('D25').Value = '2.125.57'
$ws.Range('E25').Value = '  +0.22%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '156.46'
$ws.Range('E26').Value = '  -0.12%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.376'
$ws.Range('E27').Value = '  -5.48%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.97'
$ws.Range('E28').Value = '  -0.70%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.083'
$ws.Range('E29').Value = '  -2.02%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '120.72'
$ws.Range('E30').Value = '  +0.04%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.003'
$ws.Range('E31').Value = '  -2.06%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09592'
$ws.Range('E32').Value = '  +0.31%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.602'
$ws.Range('E33').Value = '  +1.06%  '

$ws.Range('E34').Value = '  -0.35%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.390'
$ws.Range('E35').Value = '  -0.47%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06494'
$ws.Range('E36').Value = '  +5.66%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02279'
$ws.Range('E37').Value = '  -0.43%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.212'

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5919'

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '10.73'
$ws.Range('E40').Value = '  -1.17%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.852'
$ws.Range('E41').Value = '  -1.96%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.523'
$ws.Range('E42').Value = '  +3.48%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1839'
$ws.Range('E43').Value = '  -1.29%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.246'
$ws.Range('E44').Value = '  -2.76%  '

$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '12.37'
$ws.Range('E45').Value = '  -0.43%  '

$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.07530'
$ws.Range('E46').Value = '  -1.24%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5539'
$ws.Range('E47').Value = '  -0.96%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.955'
$ws.Range('E48').Value = '  -0.43%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '117.39'
$ws.Range('E49').Value = '  +0.45%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.424'
$ws.Range('E50').Value = '  -0.32%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '71.78'
$ws.Range('E51').Value = '  -1.43%  '
